$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table gained a new "2021" column (R), following the same pattern as
# the existing 2020 column (Q): a year in row 4 and a percentage in row 5.
# Copy the formatting from column Q so the new column matches visually,
# then overwrite the copied values with the new data.
$ws.Range("Q4").Copy($ws.Range("R4"))
$ws.Range("R4").Value = 2021

$ws.Range("Q5").Copy($ws.Range("R5"))
$ws.Range("R5").Value = 20.5

# Move the active selection to reflect the updated view state.
$ws.Range("S12").Select()
